# Remove the "Shoalhaven Starches Flour Mill" entry from the
# "BIO - Flour mill" sheet (row 2), and shift everything below it up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BIO - Flour mill")
$ws.Activate()

# The rows below the one being removed carry hyperlinks in column J.
# Detach them first so we can reattach them at their new (shifted-up)
# locations once the row is gone.
$ws.Range("J3").Hyperlinks.Delete()
$ws.Range("J4").Hyperlinks.Delete()
$ws.Range("J5").Hyperlinks.Delete()

# Delete the whole row - everything beneath shifts up by one.
$ws.Rows("2:2").Delete()

# Re-create the three hyperlinks one row higher than before, pointing at
# the same targets as previously.
$ws.Hyperlinks.Add($ws.Range("J2"), "https://www.inside.beer/news/detail/australia/scotland-united-malt-group-to-add-130000-tons-at-three-sites", ":~:text=3rd%2C%202020-,Australia/Scotland:%20United%20Malt%20Group%20to%20add%20130%2C000%20tons%20at,and%20in%20exports%20to%20Asia. ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J3"), "https://www.wholegrain.com.au/", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J4"), "https://www.pc.gov.au/inquiries/completed/tasmanian-shipping/submissions/submissions-test2/submission-counter/subdr071-tasmanian-shipping.docx", ":~:text=Tasmanian%20Flour%20Mills%20is%20the,are%20supplied%20with%20bagged%20flour.") | Out-Null

# Keep the same "Hyperlink" look the cells had before (Hyperlinks.Add
# above re-applies it anyway, this just normalises the style).
$ws.Range("J2:J4").Style = "Hyperlink"

# Leave the cursor where the author left it after the edit.
$ws.Range("I15").Select()
